$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.305135011672974
$ws.Range("B1").Value = 1.860609531402588
$ws.Range("C1").Value = 2.641044855117798
$ws.Range("D1").Value = 4.885111808776855
$ws.Range("E1").Value = 1.161347508430481
